$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The sheet currently named "总计" (summary) becomes the new "2022-Q1"
#    per-fund holdings sheet; a brand-new "总计" sheet is created after it
#    with the refreshed summary table (including the new 2022-Q1 row).
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$fundData = @(
    @("501054", "东方红睿泽三年定期开放灵活配置混合A", "109.00", "95.90", "6.03", "6.5727", 1),
    @("009576", "东方红智远三年持有期混合", "66.98", "92.53", "5.95", "3.9853", 1),
    @("169104", "东方红睿满沪港深灵活配置混合（LOF）", "48.91", "92.40", "6.38", "3.1205", 2),
    @("450002", "国富弹性市值混合", "41.21", "89.44", "5.29", "2.1800", 4),
    @("450009", "国富中小盘股票", "39.73", "89.63", "4.89", "1.9428", 7),
    @("310308", "申万菱信盛利精选混合", "14.72", "66.62", "4.83", "0.7110", 3),
    @("011152", "富兰克林国海兴海回报混合", "17.18", "84.91", "3.76", "0.6460", 9),
    @("910006", "东方红启盛三年持有期混合型证券投资基金A", "14.02", "91.08", "4.53", "0.6351", 5),
    @("011468", "富兰克林国海竞争优势三年持有期混合型证券投资基金A", "13.11", "85.53", "2.99", "0.3920", 9),
    @("010442", "东方红启盛三年持有期混合型证券投资基金B", "6.36", "91.08", "4.53", "0.2881", 5),
    @("163801", "中银中国混合(LOF)", "10.14", "89.19", "2.72", "0.2758", 9),
    @("910024", "东方红启阳三年持有期混合A", "6.02", "91.72", "4.58", "0.2757", 5),
    @("001726", "汇添富新兴消费股票", "3.58", "92.89", "5.11", "0.1829", 4),
    @("550001", "信诚四季红混合", "5.00", "72.84", "2.83", "0.1415", 6),
    @("487021", "工银瑞信优质精选混合", "5.96", "76.00", "1.98", "0.1180", 9),
    @("000763", "工银新财富灵活配置混合", "2.96", "92.68", "3.65", "0.1080", 6),
    @("004769", "申万菱信价值优先混合", "1.96", "91.57", "1.96", "0.0384", 10),
    @("011284", "中信保诚龙腾精选混合", "1.22", "75.38", "2.83", "0.0345", 6),
    @("006209", "中信保诚新蓝筹灵活配置混合", "1.16", "77.03", "2.89", "0.0335", 6),
    @("011032", "东方红睿泽三年定期开放灵活配置混合C", "0.35", "95.90", "6.03", "0.0211", 1),
    @("011469", "富兰克林国海竞争优势三年持有期混合型证券投资基金C", "0.70", "85.53", "2.99", "0.0209", 9),
    @("501039", "汇添富睿丰混合（LOF）A", "0.28", "20.67", "2.98", "0.0083", 3),
    @("501040", "汇添富睿丰混合（LOF）C", "0.17", "20.67", "2.98", "0.0051", 3),
    @("005247", "国都量化精选混合", "0.02", "64.74", "2.95", "0.0006", 4),
    @("003684", "汇安丰融灵活配置混合A", "0.01", "77.88", "2.91", "0.0003", 10),
    @("003685", "汇安丰融灵活配置混合C", "0.00", "77.88", "2.91", "NUM0", 10),
    @("010862", "东方红启阳三年持有期混合B", $null, "91.72", "4.58", "NUM0", 5),
)

# Header row (B1:H1), bold + centered + thin border, matching the other
# quarterly sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q1.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Pre-format the text-typed columns (B fund code, D/E/F/G numeric-looking
# text) as Text so Excel does not silently coerce "009576" -> 9576 or
# "109.00" -> 109.
$lastRow = $fundData.Length + 1
$q1.Range("B2:B$lastRow").NumberFormat = "@"
$q1.Range("D2:G$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $r = $i + 2
    $row = $fundData[$i]

    # A: zero-based row index, bold/centered/bordered like the header
    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $q1.Cells.Item($r, 2).Value = $row[0]   # B: fund code
    $q1.Cells.Item($r, 3).Value = $row[1]   # C: fund name

    if ($null -ne $row[2]) {
        $q1.Cells.Item($r, 4).Value = $row[2]   # D: fund size (text)
    }

    $q1.Cells.Item($r, 5).Value = $row[3]   # E: total stock position (text)
    $q1.Cells.Item($r, 6).Value = $row[4]   # F: position pct (text)

    if ($row[5] -eq "NUM0") {
        $q1.Cells.Item($r, 7).NumberFormat = "General"
        $q1.Cells.Item($r, 7).Value = 0     # G: held value, numeric 0
    } else {
        $q1.Cells.Item($r, 7).Value = $row[5]   # G: held value (text)
    }

    $q1.Cells.Item($r, 8).Value = $row[6]   # H: position rank (number)
}

# ---------------------------------------------------------------------------
# 2. Create the brand-new "总计" (summary) sheet right after "2022-Q1" and
#    populate it with the historical counts plus the new 2022-Q1 entry.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
for ($c = 2; $c -le 4; $c++) {
    $cell = $total.Cells.Item(1, $c)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$summaryData = @(
    @("2022-Q1", 27, 21.74),
    @("2021-Q4", 45, 57.6),
    @("2021-Q3", 96, 110.42),
    @("2021-Q2", 109, 151.45),
    @("2021-Q1", 107, 94.54000000000001),
    @("2020-Q4", 74, 59.95),
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $r = $i + 2
    $row = $summaryData[$i]

    $idxCell = $total.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}
